# Excel COM-interop script: read data from excel sheet and pass to method
# Fixes the sheet name typos ("creads" -> "creds") and tidies up the
# worksheet views / leftover formatting that Excel re-saved along with it.

$wb = $excel.ActiveWorkbook

$wsValid = $wb.Worksheets.Item(1)
$wsInvalid = $wb.Worksheets.Item(2)

# Correct the misspelled sheet names.
$wsValid.Name = "validcreds"
$wsInvalid.Name = "invalidcreds"

# Clear the stray fill formatting Excel had applied to the last three data
# rows of the "invalidcreds" sheet (A6:B8) - this merges their style back
# with the plain bordered/text-format style used by the rest of the sheet.
$wsInvalid.Range("A6:B8").Interior.Pattern = -4142  # xlNone

# Re-create the window/view state: the user had been looking at the
# "invalidcreds" sheet (zoomed in, cell E29 selected) before finally
# switching back to and leaving "validcreds" active with D28 selected.
[void]$wsInvalid.Activate()
$excel.ActiveWindow.Zoom = 142
$wsInvalid.Range("E29").Select() | Out-Null

[void]$wsValid.Activate()
$wsValid.Range("D28").Select() | Out-Null
